# Remove the "Shoes" and "Shoes Vitals" worksheets, leaving only "Shoes Gear".
$wb = $excel.ActiveWorkbook

# Deleting a sheet normally raises a confirmation prompt in Excel; suppress it.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Shoes").Delete()
$wb.Worksheets.Item("Shoes Vitals").Delete()

$excel.DisplayAlerts = $true
